# Update "想去人数" (interest count) figures in column F for both the
# "展览" (exhibition-only) sheet and the "全部类型" (all-types) sheet.
# Row numbering differs by one offset between the two sheets because the
# "全部类型" sheet contains an extra performance-event row ("演出") that
# the "展览" sheet does not.

$wb = $excel.ActiveWorkbook

$ws_exhibit = $wb.Worksheets.Item("展览")
$ws_all = $wb.Worksheets.Item("全部类型")

# Row (on "展览" sheet) -> new F value
$exhibitUpdates = @{
    2  = 168
    5  = 1832
    9  = 2413
    10 = 139
    12 = 159
    13 = 1462
    14 = 514
    17 = 223
    20 = 200
    21 = 213
    22 = 212
    24 = 114
    26 = 1509
    29 = 309
    30 = 185
    32 = 377
}

foreach ($row in $exhibitUpdates.Keys) {
    $ws_exhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# Row (on "全部类型" sheet) -> new F value
$allUpdates = @{
    2  = 168
    5  = 1832
    10 = 2413
    11 = 139
    13 = 159
    14 = 1462
    15 = 514
    18 = 223
    21 = 200
    22 = 213
    23 = 212
    25 = 114
    27 = 1509
    30 = 309
    31 = 185
    33 = 377
}

foreach ($row in $allUpdates.Keys) {
    $ws_all.Cells.Item($row, 6).Value = $allUpdates[$row]
}
